# Generate Report for Handback
#
# The localization status report is regenerated after a successful
# handback: the "Status" column flips from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, the per-language
# "Latest Handback DateTime" stamps advance to the new handback run,
# any stale "Error Detail" message is cleared now that the handback
# file is in sync, and the widened "Status" / "Error Detail" columns
# are resized to fit the new text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: both the zh-cn and de-de status columns (E, F)
# move from "Ready for handoff" to the new handed-back status.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------
# zh-cn detail sheet.
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-09-05 20:58:17"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1").ColumnWidth = 13.7470528738839

# ---------------------------------------------------------------
# de-de detail sheet.
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-09-05 20:58:25"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1").ColumnWidth = 13.7470528738839

Write-Output "Generated handback report updates."
